$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in column C (dates stored as serial numbers)
$ws.Range("C1").Value = 43936
$ws.Range("C3").Value = 43961
$ws.Range("C5").Value = 43910

# Move the active selection to C6
$ws.Range("C6").Select()
